$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-10
# from 45207 (2023-10-08) to 45208 (2023-10-09)
$ws.Range("C2:C10").Value = 45208
